# Rename "currency_movements" sheet to "money_transfers", and make it the
# active/selected sheet (it was previously "espp" that was active/selected).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("currency_movements")
$ws.Name = "money_transfers"

$ws.Activate()
$ws.Select()
